$d = $word.ActiveDocument

# Locate "Grafik och UI" (the SPRINT 3 line) in the document.
$search = $d.Content
$found = $search.Find.Execute("Grafik och UI", $true, $false, $false, $false, $false,
                               $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target text 'Grafik och UI'"
}

# Re-materialize the match bounds as a fresh Range so InsertXML replaces
# exactly that span (a Range still attached to a completed Find can behave
# like an insertion point instead of a replace-range).
$target = $d.Range($search.Start, $search.End)

# Replace the single run with four runs carrying the new wording:
#   "Grafik" + "," + " UI" + " och möjligtvis AI"
# (matches the target OOXML: same paragraph, new run split, no rPr change)
$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:r><w:t>Grafik</w:t></w:r>" +
       "<w:r><w:t>,</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> UI</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> och m&#246;jligtvis AI</w:t></w:r>" +
       "</w:p>"

$target.InsertXML($xml)
